# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# D column: price text (apostrophe-prefixed so Excel keeps exact string
# formatting like trailing zeros / many decimals instead of auto-numifying),
# then Style reset to "Normal" so no explicit cell style/quote-prefix xf is
# persisted (matches the source workbook's un-styled data cells).
# E column: volume/% text, always non-numeric-looking so plain .Value is safe.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.411.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +6.08%  '
$ws.Range("D3").Value = "'2.379.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.07%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'113.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.21%  '
$ws.Range("D6").Value = "'317.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("D7").Value = "'0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +4.06%  '
$ws.Range("D10").Value = "'42.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.94%  '
$ws.Range("D11").Value = "'0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.53%  '
$ws.Range("D12").Value = "'8.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.49%  '
$ws.Range("E13").Value = '  +5.05%  '
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = "'15.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.61%  '
$ws.Range("D16").Value = "'2.740.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.18%  '
$ws.Range("D17").Value = "'2.369.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.02%  '
$ws.Range("D18").Value = "'45.292.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.21%  '
$ws.Range("D19").Value = "'7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.43%  '
$ws.Range("D20").Value = "'0.0000108"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.48%  '
$ws.Range("D21").Value = "'13.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").Value = "'74.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.53%  '
$ws.Range("E23").Value = '  +4.26%  '
$ws.Range("D24").Value = "'269.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.23%  '
$ws.Range("E25").Value = '  +8.74%  '
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").Value = "'7.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.79%  '
$ws.Range("D28").Value = "'11.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.87%  '
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("D30").Value = "'39.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.18%  '
$ws.Range("D31").Value = "'22.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.34%  '
$ws.Range("D32").Value = "'0.0960"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.13%  '
$ws.Range("D33").Value = "'171.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.43%  '
$ws.Range("D34").Value = "'2.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +16.14%  '
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("E36").Value = '  +8.62%  '
$ws.Range("E37").Value = '  +10.19%  '
$ws.Range("E38").Value = '  +11.95%  '
$ws.Range("D39").Value = "'4.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.51%  '
$ws.Range("E40").Value = '  +5.87%  '
$ws.Range("E41").Value = '  +10.30%  '
$ws.Range("D42").Value = "'104.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.45%  '
$ws.Range("E43").Value = '  +6.76%  '
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("D45").Value = "'13.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.79%  '
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Value = "'5.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +13.53%  '
$ws.Range("D48").Value = "'116.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.81%  '
$ws.Range("D49").Value = "'1.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +19.09%  '
$ws.Range("D50").Value = "'9.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.30%  '
$ws.Range("D51").Value = "'79.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.37%  '
